$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.715.18'
$ws.Range("E2").Value = '  +2.89%  '

$ws.Range("D3").Value = '2.215.52'
$ws.Range("E3").Value = '  +0.66%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.14'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.78%  '

$ws.Range("E6").Value = '  +0.30%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '74.61'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.39%  '

$ws.Range("E8").Value = '  -0.18%  '

$ws.Range("E9").Value = '  +3.11%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.31'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.10%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0928'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.16%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.54'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.71%  '

$ws.Range("E13").Value = '  +1.05%  '

$ws.Range("E14").Value = '  -1.73%  '

$ws.Range("D15").Value = '2.545.17'
$ws.Range("E15").Value = '  +0.61%  '

$ws.Range("E16").Value = '  +4.36%  '

$ws.Range("D17").Value = '2.216.99'
$ws.Range("E17").Value = '  +1.00%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.803'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.23%  '

$ws.Range("D19").Value = '42.560.66'
$ws.Range("E19").Value = '  +2.77%  '

$ws.Range("E20").Value = '  +0.93%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '70.81'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.60%  '

$ws.Range("E22").Value = '  -2.25%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.79'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.98%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '229.71'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.99%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.19'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +7.80%  '

$ws.Range("E26").Value = '  -0.03%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.57%  '

$ws.Range("E28").Value = '  -6.95%  '

$ws.Range("E29").Value = '  -0.98%  '

$ws.Range("E30").Value = '  -0.90%  '

$ws.Range("E31").Value = '  +3.97%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '36.81'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +19.59%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.29'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.06%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0798'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.38%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.29'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.24%  '

$ws.Range("E36").Value = '  -0.66%  '

$ws.Range("E37").Value = '  +1.62%  '

$ws.Range("E38").Value = '  +5.26%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0324'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +9.22%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.77'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.24%  '

$ws.Range("E41").Value = '  +2.03%  '

$ws.Range("E42").Value = '  -0.06%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '60.79'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.69%  '

$ws.Range("E44").Value = '  +2.44%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.63'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.15%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0991'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.35%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '98.98'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.78%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.10'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.57%  '

$ws.Range("E49").Value = '  +0.10%  '

$ws.Range("E50").Value = '  -0.90%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.432'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +19.67%  '
